$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 19900
$ws.Range("J7").Value = 19900
$ws.Range("L7").Value = 19900
$ws.Range("N7").Value = -20124
$ws.Range("H12").Value = 1360.0834
$ws.Range("I12").Value = 931.5714
$ws.Range("J12").Value = 1960
$ws.Range("K12").Value = 931.5714
$ws.Range("L12").Value = 1960
$ws.Range("M12").Value = -761.5714
$ws.Range("N12").Value = -2300
$ws.Range("H14").Value = 19900
$ws.Range("J14").Value = 19900
$ws.Range("L14").Value = 19900
$ws.Range("N14").Value = -20282
$ws.Range("H55").Value = 612.25
$ws.Range("I55").Value = 225
$ws.Range("J55").Value = 999.5
$ws.Range("K55").Value = 225
$ws.Range("L55").Value = 999.5
$ws.Range("M55").Value = -11
$ws.Range("N55").Value = -1427.5
$ws.Range("H137").Value = 1856.0975
$ws.Range("I137").Value = 1305.1111
$ws.Range("J137").Value = 2287.3044
$ws.Range("K137").Value = 3915.3333
$ws.Range("L137").Value = 6861.9132
$ws.Range("M137").Value = -1365.3333
$ws.Range("N137").Value = -11961.9132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 22333.334
$ws.Range("I36").Value = 22333.334
$ws.Range("K36").Value = 22333.334
$ws.Range("M36").Value = -21987.334
$ws.Range("H45").Value = 2140
$ws.Range("I45").Value = 1936.6666
$ws.Range("J45").Value = 2750
$ws.Range("K45").Value = 1936.6666
$ws.Range("L45").Value = 2750
$ws.Range("M45").Value = -1559.6666
$ws.Range("N45").Value = -3504
$ws.Range("H61").Value = 2566.3845
$ws.Range("I61").Value = 3044.2666
$ws.Range("J61").Value = 1914.7273
$ws.Range("K61").Value = 3044.2666
$ws.Range("L61").Value = 1914.7273
$ws.Range("M61").Value = -2832.2666
$ws.Range("N61").Value = -2338.7273
$ws.Range("H74").Value = 1437.4746
$ws.Range("I74").Value = 1304.9767
$ws.Range("J74").Value = 1793.5625
$ws.Range("K74").Value = 1304.9767
$ws.Range("L74").Value = 1793.5625
$ws.Range("M74").Value = -430.9766999999999
$ws.Range("N74").Value = -3541.5625
$ws.Range("H77").Value = 1437.4746
$ws.Range("I77").Value = 1304.9767
$ws.Range("J77").Value = 1793.5625
$ws.Range("K77").Value = 6524.8835
$ws.Range("L77").Value = 8967.8125
$ws.Range("M77").Value = -2156.8835
$ws.Range("N77").Value = -17703.8125
$ws.Range("H136").Value = 2566.3845
$ws.Range("I136").Value = 3044.2666
$ws.Range("J136").Value = 1914.7273
$ws.Range("K136").Value = 9132.799800000001
$ws.Range("L136").Value = 5744.1819
$ws.Range("M136").Value = -6582.799800000001
$ws.Range("N136").Value = -10844.1819
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 16473.139
$ws.Range("I107").Value = 20390.678
$ws.Range("J107").Value = 2761.75
$ws.Range("K107").Value = 20390.678
$ws.Range("L107").Value = 2761.75
$ws.Range("M107").Value = -18470.678
$ws.Range("N107").Value = -6601.75
$ws.Range("H134").Value = 3046
$ws.Range("I134").Value = 2823.6
$ws.Range("J134").Value = 3416.6667
$ws.Range("K134").Value = 8470.799999999999
$ws.Range("L134").Value = 10250.0001
$ws.Range("M134").Value = -5935.799999999999
$ws.Range("N134").Value = -15320.0001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2230.6394
$ws.Range("J31").Value = 2052.25
$ws.Range("L31").Value = 2052.25
$ws.Range("N31").Value = -2642.25
$ws.Range("H34").Value = 2230.6394
$ws.Range("J34").Value = 2052.25
$ws.Range("L34").Value = 2052.25
$ws.Range("N34").Value = -2456.25
$ws.Range("H58").Value = 1736.7826
$ws.Range("I58").Value = 1832
$ws.Range("J58").Value = 1519.1428
$ws.Range("K58").Value = 1832
$ws.Range("L58").Value = 1519.1428
$ws.Range("M58").Value = -1629
$ws.Range("N58").Value = -1925.1428
$ws.Range("H99").Value = 1982.4
$ws.Range("I99").Value = 1982.4
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1982.4
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -484.4000000000001
$ws.Range("H126").Value = 1982.4
$ws.Range("I126").Value = 1982.4
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5947.200000000001
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -3477.200000000001
$ws.Range("H132").Value = 4006.0908
$ws.Range("I132").Value = 3302.7693
$ws.Range("K132").Value = 9908.3079
$ws.Range("M132").Value = -7378.3079
$ws.Range("H134").Value = 2202.25
$ws.Range("I134").Value = 1896.7646
$ws.Range("J134").Value = 3933.3333
$ws.Range("K134").Value = 5690.293799999999
$ws.Range("L134").Value = 11799.9999
$ws.Range("M134").Value = -3155.293799999999
$ws.Range("N134").Value = -16869.9999
$ws.Range("H136").Value = 1736.7826
$ws.Range("I136").Value = 1832
$ws.Range("J136").Value = 1519.1428
$ws.Range("K136").Value = 5496
$ws.Range("L136").Value = 4557.428400000001
$ws.Range("M136").Value = -2946
$ws.Range("N136").Value = -9657.428400000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 4000
$ws.Range("J48").Value = 4000
$ws.Range("L48").Value = 12000
$ws.Range("N48").Value = -12500
$ws.Range("H64").Value = 6000
$ws.Range("I64").Value = 2000
$ws.Range("J64").Value = 8000
$ws.Range("K64").Value = 6000
$ws.Range("L64").Value = 24000
$ws.Range("M64").Value = -5730
$ws.Range("N64").Value = -24540
$ws.Range("H67").Value = 6000
$ws.Range("I67").Value = 2000
$ws.Range("J67").Value = 8000
$ws.Range("K67").Value = 6000
$ws.Range("L67").Value = 24000
$ws.Range("M67").Value = -5064
$ws.Range("N67").Value = -25872
$ws.Range("H96").Value = 6516
$ws.Range("I96").Value = 50000
$ws.Range("J96").Value = 4625.391
$ws.Range("K96").Value = 150000
$ws.Range("L96").Value = 13876.173
$ws.Range("M96").Value = -147941
$ws.Range("N96").Value = -17994.173
$ws.Range("H113").Value = 208980
$ws.Range("I113").Value = 345397.12
$ws.Range("J113").Value = 764.3684
$ws.Range("K113").Value = 1036191.36
$ws.Range("L113").Value = 2293.1052
$ws.Range("M113").Value = -1034021.36
$ws.Range("N113").Value = -6633.1052
$ws.Range("H119").Value = 5042.7144
$ws.Range("I119").Value = 2259.8
$ws.Range("K119").Value = 6779.400000000001
$ws.Range("M119").Value = -1941.400000000001
$ws.Range("H131").Value = 21518.307
$ws.Range("J131").Value = 28846.445
$ws.Range("L131").Value = 86539.33499999999
$ws.Range("N131").Value = -96619.33499999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 50000
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H41").Value = 3648.8
$ws.Range("I41").Value = 2081.3333
$ws.Range("J41").Value = 6000
$ws.Range("K41").Value = 2081.3333
$ws.Range("L41").Value = 6000
$ws.Range("M41").Value = -1726.3333
$ws.Range("N41").Value = -6710
$ws.Range("H51").Value = 42408.668
$ws.Range("J51").Value = 42408.668
$ws.Range("L51").Value = 42408.668
$ws.Range("N51").Value = -43426.668
$ws.Range("H102").Value = 3428.476
$ws.Range("I102").Value = 3323.4119
$ws.Range("K102").Value = 3323.4119
$ws.Range("M102").Value = -1701.4119
$ws.Range("H123").Value = 8823.473
$ws.Range("J123").Value = 8823.473
$ws.Range("L123").Value = 8823.473
$ws.Range("N123").Value = -13723.473
$ws.Range("H126").Value = 4439
$ws.Range("I126").Value = 1499.5
$ws.Range("J126").Value = 5418.8335
$ws.Range("K126").Value = 4498.5
$ws.Range("L126").Value = 16256.5005
$ws.Range("M126").Value = -2028.5
$ws.Range("N126").Value = -21196.5005
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6073
$ws.Range("J7").Value = 7000.8335
$ws.Range("L7").Value = 7000.8335
$ws.Range("N7").Value = -7224.8335
$ws.Range("H22").Value = 1346.8
$ws.Range("I22").Value = 1280.2
$ws.Range("J22").Value = 1380.1
$ws.Range("K22").Value = 1280.2
$ws.Range("L22").Value = 1380.1
$ws.Range("M22").Value = -985.2
$ws.Range("N22").Value = -1970.1
$ws.Range("H27").Value = 1346.8
$ws.Range("I27").Value = 1280.2
$ws.Range("J27").Value = 1380.1
$ws.Range("K27").Value = 1280.2
$ws.Range("L27").Value = 1380.1
$ws.Range("M27").Value = -1173.2
$ws.Range("N27").Value = -1594.1
$ws.Range("H40").Value = 5000
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H126").Value = 6073
$ws.Range("J126").Value = 7000.8335
$ws.Range("L126").Value = 21002.5005
$ws.Range("N126").Value = -25942.5005
$ws.Range("H132").Value = 4416.206
$ws.Range("I132").Value = 4239.793
$ws.Range("K132").Value = 12719.379
$ws.Range("M132").Value = -10189.379
$ws.Range("H136").Value = 1592.1428
$ws.Range("I136").Value = 1406.9231
$ws.Range("K136").Value = 4220.7693
$ws.Range("M136").Value = -1670.7693
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 10387.2
$ws.Range("I126").Value = 11846.77
$ws.Range("K126").Value = 35540.31
$ws.Range("M126").Value = -33070.31
$ws.Range("H132").Value = 2885.3447
$ws.Range("I132").Value = 2326.0908
$ws.Range("J132").Value = 4643
$ws.Range("K132").Value = 6978.2724
$ws.Range("L132").Value = 13929
$ws.Range("M132").Value = -4448.2724
$ws.Range("N132").Value = -18989
$ws.Range("H136").Value = 2171.2188
$ws.Range("I136").Value = 1915.2222
$ws.Range("K136").Value = 5745.6666
$ws.Range("M136").Value = -3195.6666
